$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered/
# centered header style already used by the other headers (copy H1's format).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns I and J for rows 2-4
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2
